$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# Add "Project List" bullet entries (column E) for six sensors that
# previously had none (or an incomplete list): Soil Moisture (row 32),
# Gas Sensor Smoke/LPG (rows 33-34), Rain Water (row 37), Pulse (row 38),
# Flex Sensor (row 40) and Water Level (row 42).
#
# Each group needs extra rows inserted below the sensor's row(s) so every
# bullet gets its own row (matching the existing layout used for the other
# sensors in the sheet). We process the anchors from bottom to top so the
# row numbers of anchors still to be handled are not disturbed by earlier
# insychronous operations.
# -----------------------------------------------------------------------

# --- Water Level (row 42): insert 4 rows below it, fill E42:E46 ---
$ws.Rows("43:46").Insert()
$ws.Cells.Item(42,5).Value2 = "1. Alarm on full-filling water tank in house"
$ws.Cells.Item(43,5).Value2 = "2. Boiler water management/automation"
$ws.Cells.Item(44,5).Value2 = "3. Home Water Storage/Tank Monitoring & Automation"
$ws.Cells.Item(45,5).Value2 = "4. River water level monitoring"
$ws.Cells.Item(46,5).Value2 = "5. Drinking Water Monitoring in Wedding Function"

# --- Flex Sensor (row 40): insert 4 rows below it, fill E40:E44 ---
$ws.Rows("41:44").Insert()
$ws.Cells.Item(40,5).Value2 = "1. Gesture controlled robot/car"
$ws.Cells.Item(41,5).Value2 = "2. Controlling 3D printed Human Hand"
$ws.Cells.Item(42,5).Value2 = "3. Virtual Reality Gaming Gloves"
$ws.Cells.Item(43,5).Value2 = "4. Knee Rehabitation Monitor"
$ws.Cells.Item(44,5).Value2 = "5. PPT/Presentation helping wearable"

# --- Pulse (row 38): insert 4 rows below it, fill E38:E42 ---
$ws.Rows("39:42").Insert()
$ws.Cells.Item(38,5).Value2 = "1. Inside wearable device"
$ws.Cells.Item(39,5).Value2 = "2. In Hospital Automation"
$ws.Cells.Item(40,5).Value2 = "3. Traking heart-beat while daily workout ( in your mobile )"
$ws.Cells.Item(41,5).Value2 = "4. controls music with your heart beat"
$ws.Cells.Item(42,5).Value2 = "5. SMS alert to relative on High Heart Rate of Grand-father"

# --- Rain Water (row 37): insert 4 rows below it, fill E37:E41 ---
$ws.Rows("38:41").Insert()
$ws.Cells.Item(37,5).Value2 = "1. Car wiper Automation"
$ws.Cells.Item(38,5).Value2 = "2. Smart drying Racks that automatically comes into house during rain"
$ws.Cells.Item(39,5).Value2 = "3. Home water Harvesting Automation"
$ws.Cells.Item(40,5).Value2 = "4. Check whether its raining or not in your mobile"
$ws.Cells.Item(41,5).Value2 = "5. Tweet for First Rain of the Season"

# --- Gas Sensor Smoke/LPG (rows 33-34): insert 3 rows below row 34, fill E33:E38 ---
$ws.Rows("35:37").Insert()
$ws.Cells.Item(33,5).Value2 = "1. Industrial and Home Security System"
$ws.Cells.Item(34,5).Value2 = "2. Challenge (bump) test for Sensor manufacturing companies"
$ws.Cells.Item(35,5).Value2 = "3. Fresh Air Monitoring / weather system"
$ws.Cells.Item(36,5).Value2 = "4. Leak Detaction in Factories"
$ws.Cells.Item(37,5).Value2 = "5. LPG leak detection in kitchen"

# --- Soil Moisture (row 32): insert 2 rows below it, fill E32:E34 ---
$ws.Rows("33:34").Insert()
$ws.Cells.Item(32,5).Value2 = "1. Plant Moisture Monitoring and Automation"
$ws.Cells.Item(33,5).Value2 = "2. Smart Irrigation System/Agriculture automation"
$ws.Cells.Item(34,5).Value2 = "3. To study Ground water recharge and Evapo-transpiration"

# -----------------------------------------------------------------------
# Update the view: scroll/selection moved as part of the author's edit.
# -----------------------------------------------------------------------
$ws.Range("E64").Select()
